# Make Programming Sexy.pptx - apply author's edit:
#  1. Refresh the cached "datetimeFigureOut" (date) and "slidenum" footer
#     field text on the slide master + every slide layout (29.09.19 -> 30.09.19,
#     'Nr.' style number placeholder -> '#' style).
#  2. Rewrite the "Objectives" slide's body: "Make Programming Sexy" becomes a
#     four-line agenda ("Model based programming approach", "NAO Robot",
#     "Java Animation", "Learning Journey (Big Picture)").
#  3. Insert a brand-new "Approach" / "Agil" slide right after it (position 4),
#     pushing the old "Journey (Stacey Matrix)" slide down to position 5.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# -- 1. Date + slide-number placeholder text on master & all layouts --------
$newDate = "30.09.19"
$newNum  = [string][char]0x2039 + "#" + [string][char]0x203A

function Update-FooterFields($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = 0
        try { $phType = $sh.PlaceholderFormat.Type } catch { $phType = 0 }
        if ($phType -eq 16) {
            # ppPlaceholderDate
            $sh.TextFrame.TextRange.Text = $newDate
        } elseif ($phType -eq 13) {
            # ppPlaceholderSlideNumber
            $sh.TextFrame.TextRange.Text = $newNum
        }
    }
}

Update-FooterFields $m.Shapes
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    Update-FooterFields $m.CustomLayouts.Item($li).Shapes
}

# -- 2. "Objectives" slide body rewrite --------------------------------------
$objectivesSlide = $p.Slides.Item(3)
$body = $objectivesSlide.Shapes.Item(2)
$body.TextFrame.TextRange.Text = "Model based programming approach`rNAO Robot`rJava Animation`rLearning Journey (Big Picture)"

# -- 3. Insert the new "Approach" slide at position 4 ------------------------
$newSlide = $p.Slides.Add(4, 2)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Approach"
$newSlide.Shapes.Item(2).TextFrame.TextRange.Text = "Agil"
